$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 10500
$ws.Cells.Item(7, 10).Value = 10500
$ws.Cells.Item(7, 12).Value = 10500
$ws.Cells.Item(7, 14).Value = -10724
$ws.Cells.Item(14, 8).Value = 10500
$ws.Cells.Item(14, 10).Value = 10500
$ws.Cells.Item(14, 12).Value = 10500
$ws.Cells.Item(14, 14).Value = -10882
$ws.Cells.Item(29, 8).Value = 250000060
$ws.Cells.Item(29, 9).Value = 250000060
$ws.Cells.Item(29, 11).Value = 750000180
$ws.Cells.Item(29, 13).Value = -749999899
$ws.Cells.Item(38, 8).Value = 9221.388999999999
$ws.Cells.Item(38, 10).Value = 5999.3335
$ws.Cells.Item(38, 12).Value = 17998.0005
$ws.Cells.Item(38, 14).Value = -18742.0005
$ws.Cells.Item(40, 8).Value = 2056
$ws.Cells.Item(40, 10).Value = 2454.3635
$ws.Cells.Item(40, 12).Value = 2454.3635
$ws.Cells.Item(40, 14).Value = -2804.3635
$ws.Cells.Item(58, 8).Value = 27382.9
$ws.Cells.Item(58, 10).Value = 45576.5
$ws.Cells.Item(58, 12).Value = 136729.5
$ws.Cells.Item(58, 14).Value = -137029.5
$ws.Cells.Item(62, 8).Value = 10433.5
$ws.Cells.Item(62, 9).Value = 5899.3335
$ws.Cells.Item(62, 11).Value = 5899.3335
$ws.Cells.Item(62, 13).Value = -5275.3335
$ws.Cells.Item(65, 8).Value = 10433.5
$ws.Cells.Item(65, 9).Value = 5899.3335
$ws.Cells.Item(65, 11).Value = 29496.6675
$ws.Cells.Item(65, 13).Value = -26376.6675
$ws.Cells.Item(98, 8).Value = 2204.75
$ws.Cells.Item(98, 9).Value = 2018.6
$ws.Cells.Item(98, 10).Value = 4997
$ws.Cells.Item(98, 11).Value = 2018.6
$ws.Cells.Item(98, 12).Value = 4997
$ws.Cells.Item(98, 13).Value = -520.5999999999999
$ws.Cells.Item(98, 14).Value = -7993
$ws.Cells.Item(100, 8).Value = 3900.5
$ws.Cells.Item(100, 9).Value = 3900.5
$ws.Cells.Item(100, 11).Value = 3900.5
$ws.Cells.Item(100, 13).Value = -3359.5
$ws.Cells.Item(111, 8).Value = 1500
$ws.Cells.Item(111, 10).Value = 1500
$ws.Cells.Item(111, 12).Value = 4500
$ws.Cells.Item(111, 14).Value = -10634
$ws.Cells.Item(113, 8).Value = 6305.647
$ws.Cells.Item(113, 10).Value = 4516.3335
$ws.Cells.Item(113, 12).Value = 4516.3335
$ws.Cells.Item(113, 14).Value = -11024.3335
$ws.Cells.Item(116, 8).Value = 5039.4
$ws.Cells.Item(116, 9).Value = 7566.3335
$ws.Cells.Item(116, 10).Value = 1249
$ws.Cells.Item(116, 11).Value = 7566.3335
$ws.Cells.Item(116, 12).Value = 1249
$ws.Cells.Item(116, 13).Value = -4124.3335
$ws.Cells.Item(116, 14).Value = -8133
$ws.Cells.Item(122, 8).Value = 2204.75
$ws.Cells.Item(122, 9).Value = 2018.6
$ws.Cells.Item(122, 10).Value = 4997
$ws.Cells.Item(122, 11).Value = 6055.799999999999
$ws.Cells.Item(122, 12).Value = 14991
$ws.Cells.Item(122, 13).Value = -3605.799999999999
$ws.Cells.Item(122, 14).Value = -19891
$ws.Cells.Item(132, 8).Value = 3966.5
$ws.Cells.Item(132, 9).Value = 3574.4
$ws.Cells.Item(132, 11).Value = 10723.2
$ws.Cells.Item(132, 13).Value = -8193.200000000001
$ws.Cells.Item(137, 8).Value = 5514.879
$ws.Cells.Item(137, 9).Value = 8624.6875
$ws.Cells.Item(137, 11).Value = 25874.0625
$ws.Cells.Item(137, 13).Value = -23324.0625
$ws.Cells.Item(138, 8).Value = 3411.1538
$ws.Cells.Item(138, 9).Value = 2356.5625
$ws.Cells.Item(138, 10).Value = 3755.5103
$ws.Cells.Item(138, 11).Value = 7069.6875
$ws.Cells.Item(138, 12).Value = 11266.5309
$ws.Cells.Item(138, 13).Value = -1929.6875
$ws.Cells.Item(138, 14).Value = -21546.5309

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3348.4443
$ws.Cells.Item(2, 9).Value = 555.5
$ws.Cells.Item(2, 10).Value = 4146.4287
$ws.Cells.Item(2, 11).Value = 555.5
$ws.Cells.Item(2, 12).Value = 4146.4287
$ws.Cells.Item(2, 13).Value = -442.5
$ws.Cells.Item(2, 14).Value = -4372.4287
$ws.Cells.Item(32, 8).Value = 4817.3447
$ws.Cells.Item(32, 9).Value = 4989
$ws.Cells.Item(32, 10).Value = 2500
$ws.Cells.Item(32, 11).Value = 4989
$ws.Cells.Item(32, 12).Value = 2500
$ws.Cells.Item(32, 13).Value = -4702
$ws.Cells.Item(32, 14).Value = -3074
$ws.Cells.Item(40, 8).Value = 25000
$ws.Cells.Item(40, 9).Value = 25000
$ws.Cells.Item(40, 11).Value = 25000
$ws.Cells.Item(40, 13).Value = -24824
$ws.Cells.Item(116, 8).Value = 3348.4443
$ws.Cells.Item(116, 9).Value = 555.5
$ws.Cells.Item(116, 10).Value = 4146.4287
$ws.Cells.Item(116, 11).Value = 555.5
$ws.Cells.Item(116, 12).Value = 4146.4287
$ws.Cells.Item(116, 13).Value = 1738.5
$ws.Cells.Item(116, 14).Value = -8734.4287
$ws.Cells.Item(122, 9).Value = 6538353.5
$ws.Cells.Item(122, 10).Value = 1649.75
$ws.Cells.Item(122, 11).Value = 19615060.5
$ws.Cells.Item(122, 12).Value = 4949.25
$ws.Cells.Item(122, 13).Value = -19612610.5
$ws.Cells.Item(122, 14).Value = -9849.25
$ws.Cells.Item(138, 8).Value = 119999.5
$ws.Cells.Item(138, 10).Value = 119999.5
$ws.Cells.Item(138, 12).Value = 119999.5
$ws.Cells.Item(138, 14).Value = -130279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3348.4443
$ws.Cells.Item(3, 9).Value = 555.5
$ws.Cells.Item(3, 10).Value = 4146.4287
$ws.Cells.Item(3, 11).Value = 555.5
$ws.Cells.Item(3, 12).Value = 4146.4287
$ws.Cells.Item(3, 13).Value = -441.5
$ws.Cells.Item(3, 14).Value = -4374.4287
$ws.Cells.Item(22, 8).Value = 757.05554
$ws.Cells.Item(22, 9).Value = 703.2727
$ws.Cells.Item(22, 11).Value = 703.2727
$ws.Cells.Item(22, 13).Value = -530.2727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1793.7778
$ws.Cells.Item(22, 9).Value = 1208.6364
$ws.Cells.Item(22, 11).Value = 1208.6364
$ws.Cells.Item(22, 13).Value = -858.6364000000001
$ws.Cells.Item(50, 8).Value = 49998.332
$ws.Cells.Item(50, 10).Value = 49998.332
$ws.Cells.Item(50, 12).Value = 49998.332
$ws.Cells.Item(50, 14).Value = -51248.332
$ws.Cells.Item(118, 8).Value = 48500
$ws.Cells.Item(118, 10).Value = 48500
$ws.Cells.Item(118, 12).Value = 48500
$ws.Cells.Item(118, 14).Value = -51814

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(70, 8).Value = 2824.75
$ws.Cells.Item(70, 9).Value = 2824.75
$ws.Cells.Item(70, 11).Value = 8474.25
$ws.Cells.Item(70, 13).Value = -8159.25
$ws.Cells.Item(73, 8).Value = 2824.75
$ws.Cells.Item(73, 9).Value = 2824.75
$ws.Cells.Item(73, 11).Value = 8474.25
$ws.Cells.Item(73, 13).Value = -7382.25
$ws.Cells.Item(74, 8).Value = 27200
$ws.Cells.Item(74, 9).Value = 14000
$ws.Cells.Item(74, 11).Value = 42000
$ws.Cells.Item(74, 13).Value = -40939
$ws.Cells.Item(75, 8).Value = 2550
$ws.Cells.Item(75, 9).Value = 100
$ws.Cells.Item(75, 11).Value = 300
$ws.Cells.Item(75, 13).Value = 698
$ws.Cells.Item(77, 8).Value = 27200
$ws.Cells.Item(77, 9).Value = 14000
$ws.Cells.Item(77, 11).Value = 126000
$ws.Cells.Item(77, 13).Value = -120696
$ws.Cells.Item(78, 8).Value = 2550
$ws.Cells.Item(78, 9).Value = 100
$ws.Cells.Item(78, 11).Value = 900
$ws.Cells.Item(78, 13).Value = 4092
$ws.Cells.Item(82, 8).Value = 22759
$ws.Cells.Item(82, 9).Value = 8665
$ws.Cells.Item(82, 10).Value = 43900
$ws.Cells.Item(82, 11).Value = 25995
$ws.Cells.Item(82, 12).Value = 131700
$ws.Cells.Item(82, 13).Value = -25589
$ws.Cells.Item(82, 14).Value = -132512
$ws.Cells.Item(85, 8).Value = 22759
$ws.Cells.Item(85, 9).Value = 8665
$ws.Cells.Item(85, 10).Value = 43900
$ws.Cells.Item(85, 11).Value = 25995
$ws.Cells.Item(85, 12).Value = 131700
$ws.Cells.Item(85, 13).Value = -24591
$ws.Cells.Item(85, 14).Value = -134508
$ws.Cells.Item(92, 8).Value = 3994.2222
$ws.Cells.Item(92, 9).Value = 2933
$ws.Cells.Item(92, 10).Value = 6116.6665
$ws.Cells.Item(92, 11).Value = 8799
$ws.Cells.Item(92, 12).Value = 18349.9995
$ws.Cells.Item(92, 13).Value = -7551
$ws.Cells.Item(92, 14).Value = -20845.9995
$ws.Cells.Item(120, 8).Value = 355638
$ws.Cells.Item(120, 9).Value = 1000059
$ws.Cells.Item(120, 11).Value = 3000177
$ws.Cells.Item(120, 13).Value = -2995339
$ws.Cells.Item(121, 8).Value = 67464
$ws.Cells.Item(121, 9).Value = 100717.1
$ws.Cells.Item(121, 11).Value = 302151.3
$ws.Cells.Item(121, 13).Value = -300841.3
$ws.Cells.Item(129, 8).Value = 3161.4707
$ws.Cells.Item(129, 10).Value = 3692.9167
$ws.Cells.Item(129, 12).Value = 11078.7501
$ws.Cells.Item(129, 14).Value = -21078.7501
$ws.Cells.Item(131, 8).Value = 3385891.8
$ws.Cells.Item(131, 10).Value = 3924483
$ws.Cells.Item(131, 12).Value = 11773449
$ws.Cells.Item(131, 14).Value = -11783529
$ws.Cells.Item(140, 8).Value = 61473.766
$ws.Cells.Item(140, 9).Value = 65065.875
$ws.Cells.Item(140, 11).Value = 195197.625
$ws.Cells.Item(140, 13).Value = -190017.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 25281.857
$ws.Cells.Item(48, 10).Value = 25281.857
$ws.Cells.Item(48, 12).Value = 25281.857
$ws.Cells.Item(48, 14).Value = -26251.857
$ws.Cells.Item(117, 8).Value = 44999.75
$ws.Cells.Item(117, 10).Value = 44999.75
$ws.Cells.Item(117, 12).Value = 44999.75
$ws.Cells.Item(117, 14).Value = -51883.75
$ws.Cells.Item(122, 8).Value = 4383.0835
$ws.Cells.Item(122, 9).Value = 4327
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 12981
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -10531
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(138, 8).Value = 78750
$ws.Cells.Item(138, 10).Value = 78750
$ws.Cells.Item(138, 12).Value = 78750
$ws.Cells.Item(138, 14).Value = -89030

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 6100.0713
$ws.Cells.Item(61, 9).Value = 2760.2
$ws.Cells.Item(61, 10).Value = 7955.5557
$ws.Cells.Item(61, 11).Value = 2760.2
$ws.Cells.Item(61, 12).Value = 7955.5557
$ws.Cells.Item(61, 13).Value = -2558.2
$ws.Cells.Item(61, 14).Value = -8359.555700000001
$ws.Cells.Item(68, 8).Value = 9211.111000000001
$ws.Cells.Item(68, 9).Value = 2900
$ws.Cells.Item(68, 10).Value = 10000
$ws.Cells.Item(68, 11).Value = 2900
$ws.Cells.Item(68, 12).Value = 10000
$ws.Cells.Item(68, 13).Value = -2151
$ws.Cells.Item(68, 14).Value = -11498
$ws.Cells.Item(71, 8).Value = 9211.111000000001
$ws.Cells.Item(71, 9).Value = 2900
$ws.Cells.Item(71, 10).Value = 10000
$ws.Cells.Item(71, 11).Value = 14500
$ws.Cells.Item(71, 12).Value = 50000
$ws.Cells.Item(71, 13).Value = -10756
$ws.Cells.Item(71, 14).Value = -57488
$ws.Cells.Item(113, 8).Value = 6100.0713
$ws.Cells.Item(113, 9).Value = 2760.2
$ws.Cells.Item(113, 10).Value = 7955.5557
$ws.Cells.Item(113, 11).Value = 2760.2
$ws.Cells.Item(113, 12).Value = 7955.5557
$ws.Cells.Item(113, 13).Value = -590.1999999999998
$ws.Cells.Item(113, 14).Value = -12295.5557

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1999
$ws.Cells.Item(81, 10).Value = 1999
$ws.Cells.Item(81, 12).Value = 3998
$ws.Cells.Item(81, 14).Value = -6120
$ws.Cells.Item(84, 8).Value = 1999
$ws.Cells.Item(84, 10).Value = 1999
$ws.Cells.Item(84, 12).Value = 19990
$ws.Cells.Item(84, 14).Value = -30598
